$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 203, shifting existing rows 203:216 down to 204:217
$ws.Rows.Item(203).Insert()

# Populate the newly inserted row 203 with the new record's data
$ws.Cells.Item(203, 1).Value = 6
$ws.Cells.Item(203, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(203, 3).Value = "Metropolitana"
$ws.Cells.Item(203, 4).Value = 44714
$ws.Cells.Item(203, 5).Value = 13
$ws.Cells.Item(203, 6).Value = 100112022
$ws.Cells.Item(203, 7).Value = "Arveja Verde"
$ws.Cells.Item(203, 8).Value = "Perfection"
$ws.Cells.Item(203, 9).Value = "Primera"
$ws.Cells.Item(203, 10).Value = 200
$ws.Cells.Item(203, 11).Value = 40000
$ws.Cells.Item(203, 12).Value = 45000
$ws.Cells.Item(203, 13).Value = 42000
$ws.Cells.Item(203, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(203, 15).Value = "Provincia de Huasco"
$ws.Cells.Item(203, 16).Value = 1680
$ws.Cells.Item(203, 17).Value = 25
$ws.Cells.Item(203, 18).Value = "Hortaliza"
